# Auto-generated edit script applying scheduled-runner market-data updates
# to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5418.913
$ws.Range("I76").Value = 4952.1333
$ws.Range("J76").Value = 6294.125
$ws.Range("K76").Value = 4952.1333
$ws.Range("L76").Value = 6294.125
$ws.Range("M76").Value = -4637.1333
$ws.Range("N76").Value = -6924.125
$ws.Range("H79").Value = 5418.913
$ws.Range("I79").Value = 4952.1333
$ws.Range("J79").Value = 6294.125
$ws.Range("K79").Value = 4952.1333
$ws.Range("L79").Value = 6294.125
$ws.Range("M79").Value = -3860.1333
$ws.Range("N79").Value = -8478.125
$ws.Range("H94").Value = 5097.5713
$ws.Range("I94").Value = 5280.5
$ws.Range("K94").Value = 5280.5
$ws.Range("M94").Value = -4829.5
$ws.Range("H112").Value = 6952.067
$ws.Range("J112").Value = 7260.4185
$ws.Range("L112").Value = 21781.2555
$ws.Range("N112").Value = -23997.2555
$ws.Range("H113").Value = 3933.3333
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 4800
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 4800
$ws.Range("M113").Value = -246
$ws.Range("N113").Value = -11308
$ws.Range("H132").Value = 6245728.5
$ws.Range("I132").Value = 6840243
$ws.Range("J132").Value = 3325
$ws.Range("K132").Value = 20520729
$ws.Range("L132").Value = 9975
$ws.Range("M132").Value = -20518199
$ws.Range("N132").Value = -15035

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1190.0625
$ws.Range("I132").Value = 1216.6207
$ws.Range("K132").Value = 3649.8621
$ws.Range("M132").Value = -1119.8621

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4221.9473
$ws.Range("I134").Value = 4270.5
$ws.Range("K134").Value = 12811.5
$ws.Range("M134").Value = -10276.5
$ws.Range("H139").Value = 120907.09
$ws.Range("J139").Value = 120907.09
$ws.Range("L139").Value = 120907.09
$ws.Range("N139").Value = -131187.09

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 516.7895
$ws.Range("I22").Value = 579.7
$ws.Range("J22").Value = 446.8889
$ws.Range("K22").Value = 579.7
$ws.Range("L22").Value = 446.8889
$ws.Range("M22").Value = -229.7
$ws.Range("N22").Value = -1146.8889
$ws.Range("H31").Value = 3707536.8
$ws.Range("I31").Value = 4350542.5
$ws.Range("K31").Value = 4350542.5
$ws.Range("M31").Value = -4350247.5
$ws.Range("H34").Value = 3707536.8
$ws.Range("I34").Value = 4350542.5
$ws.Range("K34").Value = 4350542.5
$ws.Range("M34").Value = -4350340.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 4826.6665
$ws.Range("I57").Value = 495
$ws.Range("J57").Value = 5693
$ws.Range("K57").Value = 1485
$ws.Range("L57").Value = 17079
$ws.Range("M57").Value = -926
$ws.Range("N57").Value = -18197
$ws.Range("H94").Value = 2000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 6000
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -7352
$ws.Range("H121").Value = 186371.5
$ws.Range("J121").Value = 5500
$ws.Range("L121").Value = 16500
$ws.Range("N121").Value = -19120
$ws.Range("H129").Value = 3687.4666
$ws.Range("I129").Value = 1931
$ws.Range("J129").Value = 4565.7
$ws.Range("K129").Value = 5793
$ws.Range("L129").Value = 13697.1
$ws.Range("M129").Value = -793
$ws.Range("N129").Value = -23697.1
$ws.Range("H141").Value = 3582.5833
$ws.Range("I141").Value = 3582.5833
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 10747.7499
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -5567.749899999999
$ws.Range("N141").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 9969
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 9969
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 9969
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -10193
$ws.Range("H8").Value = 9969
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 9969
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 9969
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -10247
$ws.Range("H97").Value = 1757.5172
$ws.Range("I97").Value = 1457.7368
$ws.Range("K97").Value = 1457.7368
$ws.Range("M97").Value = -961.7367999999999
$ws.Range("H132").Value = 1477.2727
$ws.Range("I132").Value = 1317.5714
$ws.Range("J132").Value = 1756.75
$ws.Range("K132").Value = 3952.7142
$ws.Range("L132").Value = 5270.25
$ws.Range("M132").Value = -1422.7142
$ws.Range("N132").Value = -10330.25
$ws.Range("H135").Value = 107838
$ws.Range("J135").Value = 107838
$ws.Range("L135").Value = 107838
$ws.Range("N135").Value = -117978
$ws.Range("H138").Value = 111111
$ws.Range("J138").Value = 111111
$ws.Range("L138").Value = 111111
$ws.Range("N138").Value = -121391
$ws.Range("H139").Value = 90697.5
$ws.Range("J139").Value = 90697.5
$ws.Range("L139").Value = 90697.5
$ws.Range("N139").Value = -100977.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 2200
$ws.Range("J3").Value = 2200
$ws.Range("L3").Value = 2200
$ws.Range("N3").Value = -2424
$ws.Range("H11").Value = 11000
$ws.Range("I11").Value = 9000
$ws.Range("J11").Value = 12333.333
$ws.Range("K11").Value = 9000
$ws.Range("L11").Value = 12333.333
$ws.Range("M11").Value = -8860
$ws.Range("N11").Value = -12613.333
$ws.Range("H15").Value = 2200
$ws.Range("J15").Value = 2200
$ws.Range("L15").Value = 2200
$ws.Range("N15").Value = -2540
$ws.Range("H22").Value = 1351.1364
$ws.Range("J22").Value = 1834.3636
$ws.Range("L22").Value = 1834.3636
$ws.Range("N22").Value = -2424.3636
$ws.Range("H27").Value = 1351.1364
$ws.Range("J27").Value = 1834.3636
$ws.Range("L27").Value = 1834.3636
$ws.Range("N27").Value = -2048.3636
$ws.Range("H82").Value = 18177.4
$ws.Range("I82").Value = 12784.25
$ws.Range("K82").Value = 12784.25
$ws.Range("M82").Value = -12423.25
$ws.Range("H85").Value = 18177.4
$ws.Range("I85").Value = 12784.25
$ws.Range("K85").Value = 12784.25
$ws.Range("M85").Value = -11536.25
$ws.Range("H132").Value = 1814.3922
$ws.Range("I132").Value = 1765.721
$ws.Range("K132").Value = 5297.163
$ws.Range("M132").Value = -2767.163
$ws.Range("H137").Value = 90000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 99999
$ws.Range("J139").Value = 99999
$ws.Range("L139").Value = 99999
$ws.Range("N139").Value = -110279

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 40000
$ws.Range("J40").Value = 40000
$ws.Range("L40").Value = 40000
$ws.Range("N40").Value = -40298
